# Auto-generated COM-interop edit script.
# Populates foresatt, barn and soknad sheets with the rows added in the commit,
# matching each cell's original value + type (number vs text) as closely as the
# Excel object model allows.

$wb = $excel.ActiveWorkbook

# Scratch cell used purely to paste a 'no special number format' cell format over
# cells we forced to Text (@) so the only lasting difference is the cell's stored
# type/value, not a lingering custom number format.
function Reset-Format($ws, $cell) {
    $blank = $ws.Cells.Item(600, 600)
    $blank.Copy()
    $cell.PasteSpecial(-4122)
}

# Sets a cell to a digit-only string value without Excel's number autodetection
# converting it to a number (and losing e.g. leading zeros).
function Set-TextValue($ws, $cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    Reset-Format $ws $cell
}

# ---- foresatt sheet ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 12
$ws.Cells.Item(2, 3).Value = '-'
$ws.Cells.Item(2, 4).Value = '-'
$ws.Cells.Item(2, 5).Value = '-'
$ws.Cells.Item(2, 6).Value = '-'
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 11
$ws.Cells.Item(3, 3).Value = 'Reidar'
$ws.Cells.Item(3, 4).Value = 'Adressesvingen 7'
Set-TextValue $ws $ws.Cells.Item(3, 5) '45454545'
Set-TextValue $ws $ws.Cells.Item(3, 6) '04048701234'
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 10
$ws.Cells.Item(4, 3).Value = 'Tormund'
$ws.Cells.Item(4, 4).Value = 'Postboks 5'
Set-TextValue $ws $ws.Cells.Item(4, 5) '48544646'
Set-TextValue $ws $ws.Cells.Item(4, 6) '20018500100'
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 9
$ws.Cells.Item(5, 3).Value = 'Finn'
$ws.Cells.Item(5, 4).Value = 'Postboks 5'
Set-TextValue $ws $ws.Cells.Item(5, 5) '47474747'
Set-TextValue $ws $ws.Cells.Item(5, 6) '10108020301'
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 8
$ws.Cells.Item(6, 3).Value = 'Filip'
$ws.Cells.Item(6, 4).Value = 'Hammeren 106'
Set-TextValue $ws $ws.Cells.Item(6, 5) '45282401'
Set-TextValue $ws $ws.Cells.Item(6, 6) '15029165444'
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 7
$ws.Cells.Item(7, 3).Value = 'Cathrine'
$ws.Cells.Item(7, 4).Value = 'Hammeren 106'
Set-TextValue $ws $ws.Cells.Item(7, 5) '95793421'
Set-TextValue $ws $ws.Cells.Item(7, 6) '12048952153'
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 6
$ws.Cells.Item(8, 3).Value = 'Janne'
$ws.Cells.Item(8, 4).Value = 'Grevlingstien 44'
Set-TextValue $ws $ws.Cells.Item(8, 5) '47432211'
Set-TextValue $ws $ws.Cells.Item(8, 6) '15119165456'
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 5
$ws.Cells.Item(9, 3).Value = 'Selma'
$ws.Cells.Item(9, 4).Value = 'Grevlingstien 44'
Set-TextValue $ws $ws.Cells.Item(9, 5) '43643522'
Set-TextValue $ws $ws.Cells.Item(9, 6) '27048932123'
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 4
$ws.Cells.Item(10, 3).Value = '-'
$ws.Cells.Item(10, 4).Value = '-'
$ws.Cells.Item(10, 5).Value = '-'
$ws.Cells.Item(10, 6).Value = '-'
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 3
$ws.Cells.Item(11, 3).Value = 'Katrine'
$ws.Cells.Item(11, 4).Value = 'Svingen 10'
Set-TextValue $ws $ws.Cells.Item(11, 5) '45452323'
Set-TextValue $ws $ws.Cells.Item(11, 6) '01019500100'
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 2
$ws.Cells.Item(12, 3).Value = 'Rebecca'
$ws.Cells.Item(12, 4).Value = 'Tiurveien 12'
$ws.Cells.Item(12, 5).Value = 99778866
$ws.Cells.Item(12, 6).Value = 30109243533
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = 'Kristoffer'
$ws.Cells.Item(13, 4).Value = 'Tiurveien 12'
$ws.Cells.Item(13, 5).Value = 99887766
$ws.Cells.Item(13, 6).Value = 15029022422
$ws.Cells.Item(1, 2).Copy()
$ws.Range("A2:A13").PasteSpecial(-4122)

# ---- barn sheet ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 6
Set-TextValue $ws $ws.Cells.Item(2, 3) '17052143210'
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = 31012344544
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 4
$ws.Cells.Item(4, 3).Value = 30062174510
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = 12122278987
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = 2032290301
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = 4042011411
$ws.Cells.Item(1, 2).Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

# ---- soknad sheet ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 6
$ws.Cells.Item(2, 3).Value = 11
$ws.Cells.Item(2, 4).Value = 12
$ws.Cells.Item(2, 5).Value = 6
$ws.Cells.Item(2, 10).Value = 'Giggles and Grins Childcare'
Set-TextValue $ws $ws.Cells.Item(2, 12) '2024-12-23'
Set-TextValue $ws $ws.Cells.Item(2, 13) '10000'
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = 9
$ws.Cells.Item(3, 4).Value = 10
$ws.Cells.Item(3, 5).Value = 5
$ws.Cells.Item(3, 10).Value = 'Giggles and Grins Childcare'
Set-TextValue $ws $ws.Cells.Item(3, 12) '2024-12-23'
$ws.Cells.Item(3, 13).Value = 10000
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 4
$ws.Cells.Item(4, 3).Value = 7
$ws.Cells.Item(4, 4).Value = 8
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(4, 10).Value = '123 Learning Center'
Set-TextValue $ws $ws.Cells.Item(4, 12) '2024-12-16'
$ws.Cells.Item(4, 13).Value = 1200000
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = 5
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 10).Value = 'Sunshine Preschool'
Set-TextValue $ws $ws.Cells.Item(5, 12) '2024-12-16'
$ws.Cells.Item(5, 13).Value = 850000
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = 3
$ws.Cells.Item(6, 4).Value = 4
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 'on'
$ws.Cells.Item(6, 7).Value = 'on'
$ws.Cells.Item(6, 10).Value = 'ABC Kindergarten'
Set-TextValue $ws $ws.Cells.Item(6, 12) '2024-12-16'
$ws.Cells.Item(6, 13).Value = 300000
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 1
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 10).Value = 'Sunshine Preschool'
Set-TextValue $ws $ws.Cells.Item(7, 12) '2024-12-16'
$ws.Cells.Item(7, 13).Value = 1000000
$ws.Cells.Item(1, 2).Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

Write-Output "edit applied"
